# Add "Q5 Master Mix" reagent info to wells F1:F24 on the "Well lookup" sheet.
# This mirrors the pattern already used for wells A1:E24 (columns D/E/I):
#   D = reagent name ("Q5 Master Mix")
#   E = concentration (ng/uL) = 11
#   I = calibration type ("AQ_SP")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Well lookup")

for ($row = 122; $row -le 145; $row++) {
    $ws.Cells.Item($row, 4).Value = "Q5 Master Mix"
    $ws.Cells.Item($row, 5).Value = 11
    $ws.Cells.Item($row, 9).Value = "AQ_SP"
}
